$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_val data (regenerated to filter save games)
# Columns: Row, B(TB), C(d2S), D(K), E(IP), G(sum)
$data = @(
    ,@(2, 3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    ,@(3, 0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 4.266567469086664)
    ,@(4, 0.3464964993005633, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1.611428759079651)
    ,@(5, 0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.034748368925986)
    ,@(6, 3.182878228561681, 9.226618575922256, 3.082599426703578, 6.48142807727062, 21.97352430845813)
    ,@(7, 0.3464964993005633, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1.611428759079651)
    ,@(8, 0.006876353814593728, 0.004309184025731883, 0.1529057820181812, 0.4998867070740569, 0.6639780269325637)
    ,@(9, 3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 22.31973251085698)
    ,@(10, 1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 10.35301142835362)
    ,@(11, 0.006876353814593728, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 7.53862211810867)
    ,@(12, 3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    ,@(13, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    ,@(14, 3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    ,@(15, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    ,@(16, 1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
    ,@(17, 1.505614041169197, 0.3375848360084654, 16.98373111632243, 6.48142807727062, 25.30835807077071)
    ,@(18, 0.02258322285507441, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1.012960547955778)
    ,@(19, 3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    ,@(20, 1.505614041169197, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 2.495991366269901)
    ,@(21, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    ,@(22, 0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 8.260465185229014)
    ,@(23, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    ,@(24, 0.1554434735375247, 0.004309184025731883, 0.1529057820181812, 0.4998867070740569, 0.8125451466554947)
    ,@(25, 0.3464964993005633, 1.65323645889881, 16.98373111632243, 6.48142807727062, 25.46489215179242)
    ,@(26, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    ,@(27, 0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
    ,@(28, 3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    ,@(29, 0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 8.260465185229014)
    ,@(30, 3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 22.31973251085698)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]   # B: TB
    $ws.Cells.Item($r, 3).Value = $entry[2]   # C: d2S
    $ws.Cells.Item($r, 4).Value = $entry[3]   # D: K
    $ws.Cells.Item($r, 5).Value = $entry[4]   # E: IP
    $ws.Cells.Item($r, 7).Value = $entry[5]   # G: sum
}
